$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(37, 'Bee Swarm Activity Acoustic Classification for an IoT-Based Farm Service', 20, 2020, 33.55299442896936)
    ,@(161, 'Tracking individual honeybees among wildflower clusters with computer vision-facilitated pollinator monitoring', 8, 2021, 22.36274373259053)
    ,@(148, 'Summer weather conditions influence winter survival of honey bees (Apis mellifera) in the northeastern United States', 8, 2021, 22.36274373259053)
    ,@(201, 'Assessing the potential for deep learning and computer vision to identify bumble bee species from images', 6, 2021, 16.9658077994429)
    ,@(150, 'The prediction of swarming in honeybee colonies using vibrational spectra', 8, 2020, 13.70619777158774)
    ,@(192, 'IoT-Based Bee Swarm Activity Acoustic Classification Using Deep Neural Networks', 6, 2021, 10.72676880222841)
    ,@(228, 'Satisfiability Logic Analysis Via Radial Basis Function Neural Network with Artificial Bee Colony Algorithm', 5, 2021, 9.020640668523679)
    ,@(274, 'Markerless tracking of an entire honey bee colony', 4, 2021, 7.314512534818943)
    ,@(401, 'Automatic Data Clustering Based Mean Best Artificial Bee Colony Algorithm', 2, 2021, 6.171935933147632)
    ,@(408, 'Detection of bee diseases with a hybrid deep learning method', 2, 2021, 6.171935933147632)
    ,@(333, 'A Novel Neural Network Training Algorithm for the Identification of Nonlinear Static Systems: Artificial Bee Colony Algorithm Based on Effective Scout Bee Stage', 3, 2021, 5.608384401114207)
    ,@(332, 'Artificial Bee Colony Optimization Algorithm Incorporated With Fuzzy Theory for Real-Time Machine Learning Control of Articulated Robotic Manipulators', 3, 2020, 5.436699164345404)
    ,@(169, 'Artificial bee colony algorithm: A component-wise analysis using diversity measurement', 7, 2020, 4.82091922005571)
    ,@(193, 'Precipitation forecast of the Wujiang River Basin based on artificial bee colony algorithm and backpropagation neural network', 6, 2020, 4.159359331476323)
    ,@(411, 'An Improved Bees Algorithm for Training Deep Recurrent Networks for Sentiment Classification', 2, 2021, 3.902256267409471)
    ,@(231, 'An Efficient Intrusion Detection Model Based on Hybridization of Artificial Bee Colony and Dragonfly Algorithms for Training Multilayer Perceptrons', 5, 2020, 3.497799442896936)
    ,@(494, 'Self-Powered Smart Beehive Monitoring and Control System (SBMaCS)', 1, 2021, 3.473467966573816)
    ,@(496, 'Audio, Image, Video, and Weather Datasets for Continuous Electronic Beehive Monitoring', 1, 2021, 3.473467966573816)
    ,@(512, 'Real-time monitoring of deformed wing virus-infected bee foraging behavior following histone deacetylase inhibitor treatment', 1, 2021, 3.473467966573816)
    ,@(519, 'An Optimized Neuro-Bee Algorithm Approach to Predict the FRP-Concrete Bond Strength of RC Beams', 1, 2022, 3.473467966573816)
    ,@(539, 'Robust Queen Bee Assisted Genetic Algorithm (QBGA) Optimized Fractional Order PID (FOPID) Controller for Not Necessarily Minimum Phase Power Converters', 1, 2021, 2.196128133704736)
    ,@(510, 'Native bees in Mediterranean semi-arid agroecosystems: Unravelling the effects of biophysical habitat, floral resource, and honeybees', 1, 2021, 2.196128133704736)
    ,@(645, 'Edge-Based Detection of Varroosis in Beehives with IoT Devices with Embedded and TPU-Accelerated Machine Learning', 0, 2021, 0.775)
    ,@(646, 'Artificial Neural Networks Hidden Unit and Weight Connection Optimization by Quasi-Refection-Based Learning Artificial Bee Colony Algorithm', 0, 2021, 0.775)
    ,@(648, 'Proposed Smart Monitoring System for the Detection of Bee Swarming', 0, 2021, 0.775)
    ,@(649, 'An Intelligent Monitoring System for Assessing Bee Hive Health', 0, 2021, 0.775)
    ,@(659, 'Acoustic Scene Classification and Visualization of Beehive Sounds Using Machine Learning Algorithms and Grad-CAM', 0, 2021, 0.775)
    ,@(664, 'AlexNet, AdaBoost and Artificial Bee Colony Based Hybrid Model for Electricity Theft Detection in Smart Grids', 0, 2022, 0.775)
    ,@(653, 'Honey Bee Colony Population Daily Loss Rate Forecasting and an Early Warning Method Using Temporal Convolutional Networks', 0, 2021, 0.49)
    ,@(657, 'RFID Technology Serving Honey Bee Research: A Comprehensive Description of a 32-Antenna System to Study Honey Bee and Queen Behavior', 0, 2021, 0.49)
    ,@(660, 'DC model for SiC MOSFETs using artificial neural network optimized by artificial bee colony algorithm', 0, 2021, 0.49)
    ,@(666, 'Intrusion Detection for Network Based on Elite Clone Artificial Bee Colony and Back Propagation Neural Network', 0, 2021, 0.49)
    ,@(669, 'Comparison of Feature Extraction Methods for Sound-Based Classification of Honey Bee Activity', 0, 2022, 0.49)
    ,@(676, 'Energy-Efficient Wireless Weight Sensor for Remote Beehive Monitoring', 0, 2021, 0.49)
    ,@(832, 'Flow-mediated olfactory communication in honeybee swarms', 0, 2021, 0.395)
    ,@(652, 'Application of artificial bee colony algorithm and back propagation neural network in color evaluation of human-machine interaction interfaces', 0, $null, 0.19)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($null -ne $row[3]) {
        $ws.Cells.Item($r, 4).Value = $row[3]
    } else {
        $ws.Cells.Item($r, 4).Value = $null
    }
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Apply formatting (style index 1: bold, border, centered/top aligned) to A37, matching column A pattern
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
